$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column 'Price' values must stay TEXT (matches source inlineStr cells),
# not be auto-coerced to numbers by Excel's input parser. Temporarily
# force Text number-format while assigning, then restore the original
# cell style so no stray formatting diff is introduced.

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.681.59'
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = '  +1.00%  '

$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.286.10'
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = '  +1.03%  '

$ws.Range("E4").Value = '  -0.24%  '

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '120.68'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  +6.25%  '

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '267.46'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  +0.97%  '

$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.648'
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = '  +4.77%  '

$ws.Range("E8").Value = '  +0.15%  '

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.631'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  +5.04%  '

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '48.48'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  +1.28%  '

$ws.Range("E11").Value = '  +2.88%  '

$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.29'
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = '  +5.80%  '

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.107'
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = '  +0.03%  '

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.61'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  +1.08%  '

$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.923'
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = '  +8.04%  '

$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.626.02'
$ws.Range("D16").Style = $origStyle

$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.275.16'
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = '  +0.57%  '

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.777.29'
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = '  +1.47%  '

$ws.Range("E19").Value = '  +3.23%  '

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.97'
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = '  -0.54%  '

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.41'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  +1.94%  '

$ws.Range("E22").Value = '  +0.83%  '

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.70'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  +2.84%  '

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.66'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  -1.02%  '

$ws.Range("E25").Value = '  +1.86%  '

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.08'
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = '  +6.61%  '

$ws.Range("E27").Value = '  +1.78%  '

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '42.76'
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = '  +3.52%  '

$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.39'
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = '  +0.00%  '

$ws.Range("E30").Value = '  +0.23%  '

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '172.72'
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  +0.47%  '

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.72'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  +2.23%  '

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0934'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  +3.46%  '

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.83'
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = '  +4.86%  '

$ws.Range("E35").Value = '  +4.24%  '

$ws.Range("E36").Value = '  +13.75%  '

$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0387'
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = '  +10.49%  '

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.62'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  +0.01%  '

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.108'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  +4.07%  '

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.58'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  +7.08%  '

$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.60'
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = '  -0.33%  '

$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.78'
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = '  -3.12%  '

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.238'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  +1.98%  '

$ws.Range("E44").Value = '  -0.12%  '

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.40'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  +2.24%  '

$ws.Range("E46").Value = '  -4.73%  '

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '76.29'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  +47.55%  '

$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.28'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  +3.27%  '

$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.56'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  -0.11%  '

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.101'
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  +1.84%  '

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '102.52'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  +2.41%  '
